# This script reproduces the diff: two new rows of data were inserted into
# column A starting at row 3 (pushing the existing column-A values down by
# two rows; columns B/C/D keep their original per-row formatting), column B
# is widened/reformatted, and the sheet's used range grows from 105 to 107
# data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of existing data rows in column A before the edit.
$lastRow = 105

# Shift column A (value + cell formatting + row height) down by two rows.
# We walk from the bottom up so that we never overwrite a source row before
# it has been read.
for ($r = $lastRow; $r -ge 1; $r--) {
    $srcCell = $ws.Cells.Item($r, 1)
    $srcValue = $srcCell.Value2
    $srcHeight = $ws.Rows.Item($r).RowHeight

    $dstRow = $r + 2
    $dstCell = $ws.Cells.Item($dstRow, 1)

    $srcCell.Copy()
    $dstCell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $dstCell.Value2 = $srcValue
    $ws.Rows.Item($dstRow).RowHeight = $srcHeight
}

$excel.CutCopyMode = $false

# Widen column B and switch it to the General number format (style 0)
# instead of the Text style it used to share with the other columns.
$ws.Columns.Item(2).ClearFormats()
$ws.Columns.Item(2).ColumnWidth = 20.1

# Drop the trailing placeholder rows that used to pad the sheet out to
# 1048576 rows; the new sheet only extends to row 107.
$ws.Rows.Item(1048576).Delete()
$ws.Rows.Item(1048575).Delete()
$ws.Rows.Item(1048574).Delete()
